# Updates Colq-Musk.xlsx LR-pair data with refreshed TPM-derived values.
# Sending/target cluster set changes: "Resolving-Mac" -> "Inflammatory-Mac",
# and "ECs" becomes a sending cluster too (previously target-only),
# extending the table from 12 data rows (2-13) to 16 data rows (2-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Colq"
$ws.Cells.Item(2, 3).Value = "Musk"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.07001433333333333
$ws.Cells.Item(2, 8).Value = 0.210043
$ws.Cells.Item(2, 9).Value = 0.170290407316124
$ws.Cells.Item(2, 10).Value = 0.170290407316124
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.3676493333333333
$ws.Cells.Item(2, 14).Value = 1.102948
$ws.Cells.Item(2, 15).Value = 0.02558020875176611
$ws.Cells.Item(2, 16).Value = 0.02558020875176611
$ws.Cells.Item(2, 17).Value = 0.02574072297377778
$ws.Cells.Item(2, 18).Value = 0.231666506764
$ws.Cells.Item(2, 19).Value = 0.004356064167569732
$ws.Cells.Item(2, 20).Value = 0.004356064167569732

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Colq"
$ws.Cells.Item(3, 3).Value = "Musk"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.07001433333333333
$ws.Cells.Item(3, 8).Value = 0.210043
$ws.Cells.Item(3, 9).Value = 0.170290407316124
$ws.Cells.Item(3, 10).Value = 0.170290407316124
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.646212666666667
$ws.Cells.Item(3, 14).Value = 10.938638
$ws.Cells.Item(3, 15).Value = 0.2536952272455287
$ws.Cells.Item(3, 16).Value = 0.2536952272455287
$ws.Cells.Item(3, 17).Value = 0.2552871490482223
$ws.Cells.Item(3, 18).Value = 2.297584341434
$ws.Cells.Item(3, 19).Value = 0.04320186358179774
$ws.Cells.Item(3, 20).Value = 0.04320186358179774

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Colq"
$ws.Cells.Item(4, 3).Value = "Musk"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.07001433333333333
$ws.Cells.Item(4, 8).Value = 0.210043
$ws.Cells.Item(4, 9).Value = 0.170290407316124
$ws.Cells.Item(4, 10).Value = 0.170290407316124
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.032567
$ws.Cells.Item(4, 14).Value = 0.097701
$ws.Cells.Item(4, 15).Value = 0.002265938172294887
$ws.Cells.Item(4, 16).Value = 0.002265938172294887
$ws.Cells.Item(4, 17).Value = 0.002280156793666666
$ws.Cells.Item(4, 18).Value = 0.020521411143
$ws.Cells.Item(4, 19).Value = 0.0003858675343132499
$ws.Cells.Item(4, 20).Value = 0.0003858675343132499

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Colq"
$ws.Cells.Item(5, 3).Value = "Musk"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.07001433333333333
$ws.Cells.Item(5, 8).Value = 0.210043
$ws.Cells.Item(5, 9).Value = 0.170290407316124
$ws.Cells.Item(5, 10).Value = 0.170290407316124
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 10.32598433333333
$ws.Cells.Item(5, 14).Value = 30.977953
$ws.Cells.Item(5, 15).Value = 0.7184586258304102
$ws.Cells.Item(5, 16).Value = 0.7184586258304102
$ws.Cells.Item(5, 17).Value = 0.7229669091087777
$ws.Cells.Item(5, 18).Value = 6.506702181979
$ws.Cells.Item(5, 19).Value = 0.1223466120324433
$ws.Cells.Item(5, 20).Value = 0.1223466120324433

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Colq"
$ws.Cells.Item(6, 3).Value = "Musk"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.1708176666666666
$ws.Cells.Item(6, 8).Value = 0.5124529999999999
$ws.Cells.Item(6, 9).Value = 0.4154665001945777
$ws.Cells.Item(6, 10).Value = 0.4154665001945778
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.3676493333333333
$ws.Cells.Item(6, 14).Value = 1.102948
$ws.Cells.Item(6, 15).Value = 0.02558020875176611
$ws.Cells.Item(6, 16).Value = 0.02558020875176611
$ws.Cells.Item(6, 17).Value = 0.06280100127155555
$ws.Cells.Item(6, 18).Value = 0.5652090114439999
$ws.Cells.Item(6, 19).Value = 0.01062771980434297
$ws.Cells.Item(6, 20).Value = 0.01062771980434298

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Colq"
$ws.Cells.Item(7, 3).Value = "Musk"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.1708176666666666
$ws.Cells.Item(7, 8).Value = 0.5124529999999999
$ws.Cells.Item(7, 9).Value = 0.4154665001945777
$ws.Cells.Item(7, 10).Value = 0.4154665001945778
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.646212666666667
$ws.Cells.Item(7, 14).Value = 10.938638
$ws.Cells.Item(7, 15).Value = 0.2536952272455287
$ws.Cells.Item(7, 16).Value = 0.2536952272455287
$ws.Cells.Item(7, 17).Value = 0.6228375398904444
$ws.Cells.Item(7, 18).Value = 5.605537859013999
$ws.Cells.Item(7, 19).Value = 0.1054018681797679
$ws.Cells.Item(7, 20).Value = 0.1054018681797679

# Row 8: FAPs -> Inflammatory-Mac
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Colq"
$ws.Cells.Item(8, 3).Value = "Musk"
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.1708176666666666
$ws.Cells.Item(8, 8).Value = 0.5124529999999999
$ws.Cells.Item(8, 9).Value = 0.4154665001945777
$ws.Cells.Item(8, 10).Value = 0.4154665001945778
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.032567
$ws.Cells.Item(8, 14).Value = 0.097701
$ws.Cells.Item(8, 15).Value = 0.002265938172294887
$ws.Cells.Item(8, 16).Value = 0.002265938172294887
$ws.Cells.Item(8, 17).Value = 0.005563018950333332
$ws.Cells.Item(8, 18).Value = 0.05006717055299999
$ws.Cells.Item(8, 19).Value = 0.0009414214021006548
$ws.Cells.Item(8, 20).Value = 0.0009414214021006549

# Row 9: FAPs -> MuSCs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Colq"
$ws.Cells.Item(9, 3).Value = "Musk"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.1708176666666666
$ws.Cells.Item(9, 8).Value = 0.5124529999999999
$ws.Cells.Item(9, 9).Value = 0.4154665001945777
$ws.Cells.Item(9, 10).Value = 0.4154665001945778
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 10.32598433333333
$ws.Cells.Item(9, 14).Value = 30.977953
$ws.Cells.Item(9, 15).Value = 0.7184586258304102
$ws.Cells.Item(9, 16).Value = 0.7184586258304102
$ws.Cells.Item(9, 17).Value = 1.763860549856555
$ws.Cells.Item(9, 18).Value = 15.874744948709
$ws.Cells.Item(9, 19).Value = 0.2984954908083661
$ws.Cells.Item(9, 20).Value = 0.2984954908083662

# Row 10: Inflammatory-Mac -> ECs
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Colq"
$ws.Cells.Item(10, 3).Value = "Musk"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.002571333333333333
$ws.Cells.Item(10, 8).Value = 0.007714
$ws.Cells.Item(10, 9).Value = 0.006254053703463484
$ws.Cells.Item(10, 10).Value = 0.006254053703463485
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.3676493333333333
$ws.Cells.Item(10, 14).Value = 1.102948
$ws.Cells.Item(10, 15).Value = 0.02558020875176611
$ws.Cells.Item(10, 16).Value = 0.02558020875176611
$ws.Cells.Item(10, 17).Value = 0.0009453489857777778
$ws.Cells.Item(10, 18).Value = 0.008508140872
$ws.Cells.Item(10, 19).Value = 0.0001599799992793519
$ws.Cells.Item(10, 20).Value = 0.0001599799992793519

# Row 11: Inflammatory-Mac -> FAPs
$ws.Cells.Item(11, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 2).Value = "Colq"
$ws.Cells.Item(11, 3).Value = "Musk"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.002571333333333333
$ws.Cells.Item(11, 8).Value = 0.007714
$ws.Cells.Item(11, 9).Value = 0.006254053703463484
$ws.Cells.Item(11, 10).Value = 0.006254053703463485
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 3.646212666666667
$ws.Cells.Item(11, 14).Value = 10.938638
$ws.Cells.Item(11, 15).Value = 0.2536952272455287
$ws.Cells.Item(11, 16).Value = 0.2536952272455287
$ws.Cells.Item(11, 17).Value = 0.009375628170222224
$ws.Cells.Item(11, 18).Value = 0.08438065353200001
$ws.Cells.Item(11, 19).Value = 0.001586623575505909
$ws.Cells.Item(11, 20).Value = 0.00158662357550591

# Row 12: Inflammatory-Mac -> Inflammatory-Mac
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "Colq"
$ws.Cells.Item(12, 3).Value = "Musk"
$ws.Cells.Item(12, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.002571333333333333
$ws.Cells.Item(12, 8).Value = 0.007714
$ws.Cells.Item(12, 9).Value = 0.006254053703463484
$ws.Cells.Item(12, 10).Value = 0.006254053703463485
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.032567
$ws.Cells.Item(12, 14).Value = 0.097701
$ws.Cells.Item(12, 15).Value = 0.002265938172294887
$ws.Cells.Item(12, 16).Value = 0.002265938172294887
$ws.Cells.Item(12, 17).Value = 0.00008374061266666666
$ws.Cells.Item(12, 18).Value = 0.000753665514
$ws.Cells.Item(12, 19).Value = 0.00001417129901826012
$ws.Cells.Item(12, 20).Value = 0.00001417129901826012

# Row 13: Inflammatory-Mac -> MuSCs
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "Colq"
$ws.Cells.Item(13, 3).Value = "Musk"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.002571333333333333
$ws.Cells.Item(13, 8).Value = 0.007714
$ws.Cells.Item(13, 9).Value = 0.006254053703463484
$ws.Cells.Item(13, 10).Value = 0.006254053703463485
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 10.32598433333333
$ws.Cells.Item(13, 14).Value = 30.977953
$ws.Cells.Item(13, 15).Value = 0.7184586258304102
$ws.Cells.Item(13, 16).Value = 0.7184586258304102
$ws.Cells.Item(13, 17).Value = 0.02655154771577778
$ws.Cells.Item(13, 18).Value = 0.238963929442
$ws.Cells.Item(13, 19).Value = 0.004493278829659963
$ws.Cells.Item(13, 20).Value = 0.004493278829659963

# Row 14: MuSCs -> ECs
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Colq"
$ws.Cells.Item(14, 3).Value = "Musk"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.1677433333333334
$ws.Cells.Item(14, 8).Value = 0.5032300000000001
$ws.Cells.Item(14, 9).Value = 0.4079890387858348
$ws.Cells.Item(14, 10).Value = 0.4079890387858348
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.3676493333333333
$ws.Cells.Item(14, 14).Value = 1.102948
$ws.Cells.Item(14, 15).Value = 0.02558020875176611
$ws.Cells.Item(14, 16).Value = 0.02558020875176611
$ws.Cells.Item(14, 17).Value = 0.06167072467111112
$ws.Cells.Item(14, 18).Value = 0.5550365220400001
$ws.Cells.Item(14, 19).Value = 0.01043644478057406
$ws.Cells.Item(14, 20).Value = 0.01043644478057406

# Row 15: MuSCs -> FAPs
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Colq"
$ws.Cells.Item(15, 3).Value = "Musk"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.1677433333333334
$ws.Cells.Item(15, 8).Value = 0.5032300000000001
$ws.Cells.Item(15, 9).Value = 0.4079890387858348
$ws.Cells.Item(15, 10).Value = 0.4079890387858348
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 3.646212666666667
$ws.Cells.Item(15, 14).Value = 10.938638
$ws.Cells.Item(15, 15).Value = 0.2536952272455287
$ws.Cells.Item(15, 16).Value = 0.2536952272455287
$ws.Cells.Item(15, 17).Value = 0.611627866748889
$ws.Cells.Item(15, 18).Value = 5.504650800740001
$ws.Cells.Item(15, 19).Value = 0.1035048719084572
$ws.Cells.Item(15, 20).Value = 0.1035048719084572

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Colq"
$ws.Cells.Item(16, 3).Value = "Musk"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.1677433333333334
$ws.Cells.Item(16, 8).Value = 0.5032300000000001
$ws.Cells.Item(16, 9).Value = 0.4079890387858348
$ws.Cells.Item(16, 10).Value = 0.4079890387858348
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.032567
$ws.Cells.Item(16, 14).Value = 0.097701
$ws.Cells.Item(16, 15).Value = 0.002265938172294887
$ws.Cells.Item(16, 16).Value = 0.002265938172294887
$ws.Cells.Item(16, 17).Value = 0.005462897136666667
$ws.Cells.Item(16, 18).Value = 0.04916607423
$ws.Cells.Item(16, 19).Value = 0.0009244779368627223
$ws.Cells.Item(16, 20).Value = 0.0009244779368627223

# Row 17: MuSCs -> MuSCs
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Colq"
$ws.Cells.Item(17, 3).Value = "Musk"
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.1677433333333334
$ws.Cells.Item(17, 8).Value = 0.5032300000000001
$ws.Cells.Item(17, 9).Value = 0.4079890387858348
$ws.Cells.Item(17, 10).Value = 0.4079890387858348
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 10.32598433333333
$ws.Cells.Item(17, 14).Value = 30.977953
$ws.Cells.Item(17, 15).Value = 0.7184586258304102
$ws.Cells.Item(17, 16).Value = 0.7184586258304102
$ws.Cells.Item(17, 17).Value = 1.732115032021111
$ws.Cells.Item(17, 18).Value = 15.58903528819
$ws.Cells.Item(17, 19).Value = 0.2931232441599408
$ws.Cells.Item(17, 20).Value = 0.2931232441599408

